# "tratamento do eixo 10"
# - Rename sheet "Sheet1" -> "2018"
# - Shorten the UF header labels:
#     A1: "Sigla da UF"        -> "Sigla"
#     B1: "Unidade da Federação" -> "UF"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "2018"

$ws.Range("A1").Value = "Sigla"
$ws.Range("B1").Value = "UF"
